$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Metadata": version bump, new publication date, publisher name, and
# replacing the duplicated "Contact" row with a single "Jurisdiction" row.
# ---------------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

# Version 5.0.0 -> 6.0.0
$meta.Cells.Item(3, 2).Value = "6.0.0"

# Date -> new publication date
$meta.Cells.Item(8, 2).Value = "2022-01-21T20:46:54+00:00"

# Publisher value was empty, now populated
$meta.Cells.Item(9, 2).Value = "Alvearie Team"

# Remove the duplicate "Contact" row (row 11); row 10 ("Contact" / "No
# display for ContactDetail") becomes the new "Jurisdiction" row after the
# delete shifts everything below it up by one.
$meta.Rows.Item(11).Delete()

$meta.Cells.Item(10, 1).Value = "Jurisdiction"
$meta.Cells.Item(10, 2).Value = "United States of America"

# ---------------------------------------------------------------------------
# Sheet "Elements": update the root Extension row's Short/Definition text to
# reflect the renamed profile.
# ---------------------------------------------------------------------------
$elements = $wb.Worksheets.Item("Elements")

$elements.Cells.Item(2, 11).Value = "Measure Population Id"
$elements.Cells.Item(2, 12).Value = "Unique static identifier for the measure population that does not change between cohort engine runs. This is a user-friendly textual key (internal use only) that identifies a specific population. This id should be unique across all measures and should not change when a measure is versioned or when overrides are applied."
